$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Selects username and password"
$rng = $d.Content
$found = $rng.Find.Execute("Selects username and password")
if (-not $found) {
    throw "Could not find the target paragraph text."
}

# Grab the whole paragraph (including its end-of-paragraph mark) so we can
# replace its run content while keeping the paragraph properties intact.
$para = $rng.Paragraphs(1).Range

# Rebuild the paragraph as two runs: "U" (keeping the original run's
# rsidRPr identity) + "sername and password" (a brand-new run), so the
# final text reads "Username and password" instead of
# "Selects username and password".
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="002B7BCC" w:rsidRDefault="002B7BCC" w:rsidP="00435742"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:after="0"/><w:ind w:left="1080"/><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00255BC7"><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>U</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>sername and password</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$para.InsertXML($xml)
